$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 158

$ws.Range("E11").Value = 447
$ws.Range("F11").Value = 253
$ws.Range("H11").Value = 318

$ws.Range("F13").Value = 93
$ws.Range("H13").Value = 127

$ws.Range("E14").Value = 142

$ws.Range("E16").Value = 232

$ws.Range("E21").Value = 150

$ws.Range("E24").Value = 260
$ws.Range("F24").Value = 153
$ws.Range("H24").Value = 183

$ws.Range("E25").Value = 325

$ws.Range("E26").Value = 187
$ws.Range("F26").Value = 118
$ws.Range("H26").Value = 143

$ws.Range("E27").Value = 378

$ws.Range("E32").Value = 213
$ws.Range("F32").Value = 135
$ws.Range("H32").Value = 173

$ws.Range("E34").Value = 246

$ws.Range("E35").Value = 179

$ws.Range("E37").Value = 191

$ws.Range("E39").Value = 196

$ws.Range("E40").Value = 303

$ws.Range("F41").Value = 218
$ws.Range("H41").Value = 310

$ws.Range("E44").Value = 359
$ws.Range("F44").Value = 190
$ws.Range("H44").Value = 258

$ws.Range("E46").Value = 385
$ws.Range("F46").Value = 225
$ws.Range("H46").Value = 289

$ws.Range("E47").Value = 532
$ws.Range("F47").Value = 300
$ws.Range("H47").Value = 392

$ws.Range("E48").Value = 262

$ws.Range("E52").Value = 32
